# "Drop in RMI script results for 3.0"
# Update the About sheet's currency-year labels and the underlying
# conversion-factor value from the 2018-dollar vintage to the 2019-dollar
# vintage. The three OCCF-Dp*OCU sheets recompute automatically because
# their cells hold formulas that reference About!A26.

$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")

# Conversion factor (2012 dollars per 2019 dollar), replacing the old
# 2012-dollars-per-2018-dollar figure.
$about.Range("A26").Value = 0.89805481563188172

# Text labels referencing the dollar vintage.
$about.Range("A18").Value = "billion 2019 dollars"
$about.Range("A21").Value = "million 2019 dollars"
$about.Range("B26").Value = "2019 dollars per 2012 dollar"
$about.Range("B29").Value = 'which in this case is "2012 dollars per 2019 dollar."'

$excel.CalculateFullRebuild()
